$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font/border/alignment) from the last header cell (L1) to the
# new header cell (M1), then set the new header text.
$ws.Range("L1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "包含公衛衛教措施比例(%)"

# New data values for column M, rows 2-18
$ws.Range("M2").Value = 6
$ws.Range("M3").Value = 6
$ws.Range("M4").Value = 7
$ws.Range("M5").Value = 7
$ws.Range("M6").Value = 9
$ws.Range("M7").Value = 11
$ws.Range("M8").Value = 13
$ws.Range("M9").Value = 11
$ws.Range("M10").Value = 10
$ws.Range("M11").Value = 25
$ws.Range("M12").Value = 30
$ws.Range("M13").Value = 17
$ws.Range("M14").Value = 22
$ws.Range("M15").Value = 22
$ws.Range("M16").Value = 23
$ws.Range("M17").Value = 24
$ws.Range("M18").Value = 34
